$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GlobalPriority (B) and NivelSeguridad (C) values to reflect the
# corrected ranking computation.

$ws.Range("B2").Value = 0.2564857595068418
$ws.Range("C2").Value = 1

$ws.Range("B3").Value = 0.158943965723677
$ws.Range("C3").Value = 2

$ws.Range("B4").Value = 0.2920350876630113
$ws.Range("C4").Value = 1

$ws.Range("B5").Value = 0.1991173787737204
$ws.Range("C5").Value = 1

$ws.Range("B6").Value = 0.29526071091122
$ws.Range("C6").Value = 1

$ws.Range("B7").Value = 0.2911332986410671
$ws.Range("C7").Value = 1

$ws.Range("B8").Value = 0.198237677006421
$ws.Range("C8").Value = 1

$ws.Range("B9").Value = 0.2347069814647371
$ws.Range("C9").Value = 1
